$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-10 Saturday" "2024-02-11 Sunday"

Replace-Text "89÷3=29, 2" "46÷9=5, 1"
Replace-Text "43÷9=4, 7" "19÷5=3, 4"
Replace-Text "35÷8=4, 3" "50÷3=16, 2"
Replace-Text "49÷6=8, 1" "97÷6=16, 1"
Replace-Text "86÷9=9, 5" "61÷6=10, 1"

Replace-Text "62÷8=7, 6" "14÷5=2, 4"
Replace-Text "21÷3=7, 0" "51÷2=25, 1"
Replace-Text "66÷3=22, 0" "34÷2=17, 0"
Replace-Text "37÷5=7, 2" "81÷6=13, 3"
Replace-Text "47÷5=9, 2" "56÷2=28, 0"

Replace-Text "57÷9=6, 3" "39÷4=9, 3"
Replace-Text "25÷9=2, 7" "50÷2=25, 0"
Replace-Text "95÷7=13, 4" "11÷9=1, 2"
Replace-Text "51÷3=17, 0" "39÷6=6, 3"
Replace-Text "29÷9=3, 2" "55÷3=18, 1"

Replace-Text "56÷7=8, 0" "89÷4=22, 1"
Replace-Text "28÷5=5, 3" "40÷3=13, 1"
Replace-Text "33÷4=8, 1" "54÷5=10, 4"
Replace-Text "80÷7=11, 3" "38÷5=7, 3"
Replace-Text "88÷3=29, 1" "34÷3=11, 1"

Replace-Text "26÷8=3, 2" "49÷3=16, 1"
Replace-Text "84÷5=16, 4" "29÷4=7, 1"
Replace-Text "23÷7=3, 2" "14÷9=1, 5"
Replace-Text "60÷7=8, 4" "51÷7=7, 2"
Replace-Text "77÷4=19, 1" "85÷8=10, 5"
